# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 07:18"

# Row 6 - India
$ws.Range("B6").Value = 3387500
$ws.Range("C6").Value = 2925
$ws.Range("D6").Value = 2583948
$ws.Range("E6").Value = 741858

# Row 19 - Pakistan
$ws.Range("B19").Value = 295053
$ws.Range("C19").Value = 415
$ws.Range("D19").Value = 279937
$ws.Range("E19").Value = 8833
$ws.Range("G19").Value = 9
$ws.Range("H19").Value = 6283

# Row 62 - Uzbekistan
$ws.Range("B62").Value = 40447
$ws.Range("D62").Value = 37200
$ws.Range("E62").Value = 2949
$ws.Range("H62").Value = 298

# Row 72 - Australia
$ws.Range("B72").Value = 25448
$ws.Range("C72").Value = 126
$ws.Range("E72").Value = 4498

# Row 123 - Tailandia
$ws.Range("B123").Value = 3410
$ws.Range("C123").Value = 6
$ws.Range("E123").Value = 115

# Row 188 - Butan
$ws.Range("B188").Value = 184
$ws.Range("C188").Value = 1
$ws.Range("D188").Value = 119
